# Apply the "respond to comments / increase coverage" edits to Sheet1 of
# the cross_sheet workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# D1 and D2 pick up the "applied" style already used by the rest of their
# rows (they previously used the default/unapplied style).
$ws.Range("D1").Style = $ws.Range("C1").Style
$ws.Range("D2").Style = $ws.Range("C2").Style

# New value in F2.
$ws.Range("F2").Value = 3

# New formula in D3 summing row 2 (depends on the new F2 value and the
# recalculated D2, so it naturally lands on 30).
$ws.Range("D3").Formula = "=SUM(2:2)"

# New value in row 7 (bumps D2's SUM(Sheet1!A:A) total up to 24).
$ws.Range("A7").Value = 16

# Update the active selection to reflect where the author ended up editing.
$ws.Range("C7").Select()
